# Applies the "tested and working with simulation" edit to
# experiment/template/experiment_configurations.xlsx
#
# Sheets:
#   selection-control        (sheet1) - per-vial step-type / growthrate config
#   selection-step_generation(sheet2) - per-step rescue settings
#   selection-steps           (sheet3) - per-step logarithmic_steps flag

$wb = $excel.ActiveWorkbook

$wsControl = $wb.Worksheets.Item("selection-control")
$wsStepGen = $wb.Worksheets.Item("selection-step_generation")
$wsSteps   = $wb.Worksheets.Item("selection-steps")

# -----------------------------------------------------------------
# Sheet "selection-control": step_type (col B) and growthrate bounds
# (cols H/I) updates.
# -----------------------------------------------------------------

# step_type (column B) changes
$wsControl.Range("B2").Value  = "AUTO"
$wsControl.Range("B3").Value  = "AUTO"
$wsControl.Range("B4").Value  = "AUTO"
$wsControl.Range("B5").Value  = "AUTO"
$wsControl.Range("B14").Value = "OFF"
$wsControl.Range("B15").Value = "OFF"
$wsControl.Range("B16").Value = "MANUAL"
$wsControl.Range("B17").Value = "MANUAL"

# min_growthrate (H) / max_growthrate (I) updated for every data row
for ($r = 2; $r -le 17; $r++) {
    $wsControl.Range("H$r").Value = 0.1
    $wsControl.Range("I$r").Value = 0.12
}

# -----------------------------------------------------------------
# Sheet "selection-step_generation": rescue_dilutions (B, bool),
# rescue_threshold (C), step_number (D), min_selection (E) updates.
# -----------------------------------------------------------------

$wsStepGen.Range("C2").Value = 1
$wsStepGen.Range("D2").Value = 20
$wsStepGen.Range("E2").Value = 20

$wsStepGen.Range("B3").Value = $true
$wsStepGen.Range("C3").Value = 1
$wsStepGen.Range("D3").Value = 50
$wsStepGen.Range("E3").Value = 20

$wsStepGen.Range("B12").Value = $true
$wsStepGen.Range("B13").Value = $true

# -----------------------------------------------------------------
# Sheet "selection-steps": logarithmic_steps (column B) updates.
# -----------------------------------------------------------------

$wsSteps.Range("B4").Value  = "-"
$wsSteps.Range("B5").Value  = "-"
$wsSteps.Range("B16").Value = "1,2,3,4"
$wsSteps.Range("B17").Value = "1,2,3,4"

# -----------------------------------------------------------------
# Selections / active sheet, applied last so the saved view state
# matches what was left selected when the workbook was saved.
# -----------------------------------------------------------------

$null = $wsControl.Activate()
$null = $wsControl.Range("F4").Select()

$null = $wsSteps.Activate()
$null = $wsSteps.Range("B3:B6").Select()

$null = $wsStepGen.Activate()
$null = $wsStepGen.Range("E3").Select()
